$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Adam9"
$ws.Range("C2").Value = "Itga3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 10.32689766666667
$ws.Range("H2").Value = 30.980693
$ws.Range("I2").Value = 0.2044815006034941
$ws.Range("J2").Value = 0.204481500603494
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 6.169512999999999
$ws.Range("N2").Value = 18.508539
$ws.Range("O2").Value = 0.5207942167525852
$ws.Range("P2").Value = 0.5207942167525853
$ws.Range("Q2").Value = 63.71192940416965
$ws.Range("R2").Value = 573.4073646375269
$ws.Range("S2").Value = 0.10649278294719
$ws.Range("T2").Value = 0.10649278294719

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Adam9"
$ws.Range("C3").Value = "Itga3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 10.32689766666667
$ws.Range("H3").Value = 30.980693
$ws.Range("I3").Value = 0.2044815006034941
$ws.Range("J3").Value = 0.204481500603494
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.06813733333333333
$ws.Range("N3").Value = 0.204412
$ws.Range("O3").Value = 0.005751755307905689
$ws.Range("P3").Value = 0.00575175530790569
$ws.Range("Q3").Value = 0.7036472686128887
$ws.Range("R3").Value = 6.332825417515999
$ws.Range("S3").Value = 0.001176127556464667
$ws.Range("T3").Value = 0.001176127556464667

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Adam9"
$ws.Range("C4").Value = "Itga3"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 10.32689766666667
$ws.Range("H4").Value = 30.980693
$ws.Range("I4").Value = 0.2044815006034941
$ws.Range("J4").Value = 0.204481500603494
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 5.608704333333333
$ws.Range("N4").Value = 16.826113
$ws.Range("O4").Value = 0.473454027939509
$ws.Range("P4").Value = 0.4734540279395091
$ws.Range("Q4").Value = 57.92051569292321
$ws.Range("R4").Value = 521.284641236309
$ws.Range("S4").Value = 0.09681259009983942
$ws.Range("T4").Value = 0.09681259009983942

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Adam9"
$ws.Range("C5").Value = "Itga3"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 29.37031933333333
$ws.Range("H5").Value = 88.110958
$ws.Range("I5").Value = 0.5815577111671272
$ws.Range("J5").Value = 0.5815577111671272
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 6.169512999999999
$ws.Range("N5").Value = 18.508539
$ws.Range("O5").Value = 0.5207942167525852
$ws.Range("P5").Value = 0.5207942167525853
$ws.Range("Q5").Value = 181.2005669411513
$ws.Range("R5").Value = 1630.805102470362
$ws.Range("S5").Value = 0.3028718926837102
$ws.Range("T5").Value = 0.3028718926837102

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Adam9"
$ws.Range("C6").Value = "Itga3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 29.37031933333333
$ws.Range("H6").Value = 88.110958
$ws.Range("I6").Value = 0.5815577111671272
$ws.Range("J6").Value = 0.5815577111671272
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.06813733333333333
$ws.Range("N6").Value = 0.204412
$ws.Range("O6").Value = 0.005751755307905689
$ws.Range("P6").Value = 0.00575175530790569
$ws.Range("Q6").Value = 2.001215238521778
$ws.Range("R6").Value = 18.010937146696
$ws.Range("S6").Value = 0.003344977652059007
$ws.Range("T6").Value = 0.003344977652059008

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Adam9"
$ws.Range("C7").Value = "Itga3"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 29.37031933333333
$ws.Range("H7").Value = 88.110958
$ws.Range("I7").Value = 0.5815577111671272
$ws.Range("J7").Value = 0.5815577111671272
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 5.608704333333333
$ws.Range("N7").Value = 16.826113
$ws.Range("O7").Value = 0.473454027939509
$ws.Range("P7").Value = 0.4734540279395091
$ws.Range("Q7").Value = 164.7294373162504
$ws.Range("R7").Value = 1482.564935846254
$ws.Range("S7").Value = 0.275340840831358
$ws.Range("T7").Value = 0.275340840831358

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Adam9"
$ws.Range("C8").Value = "Itga3"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 10.80562866666667
$ws.Range("H8").Value = 32.416886
$ws.Range("I8").Value = 0.2139607882293788
$ws.Range("J8").Value = 0.2139607882293788
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 6.169512999999999
$ws.Range("N8").Value = 18.508539
$ws.Range("O8").Value = 0.5207942167525852
$ws.Range("P8").Value = 0.5207942167525853
$ws.Range("Q8").Value = 66.66546653217266
$ws.Range("R8").Value = 599.9891987895539
$ws.Range("S8").Value = 0.1114295411216851
$ws.Range("T8").Value = 0.1114295411216851

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Adam9"
$ws.Range("C9").Value = "Itga3"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 10.80562866666667
$ws.Range("H9").Value = 32.416886
$ws.Range("I9").Value = 0.2139607882293788
$ws.Range("J9").Value = 0.2139607882293788
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.06813733333333333
$ws.Range("N9").Value = 0.204412
$ws.Range("O9").Value = 0.005751755307905689
$ws.Range("P9").Value = 0.00575175530790569
$ws.Range("Q9").Value = 0.7362667223368887
$ws.Range("R9").Value = 6.626400501031999
$ws.Range("S9").Value = 0.001230650099382015
$ws.Range("T9").Value = 0.001230650099382015

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Adam9"
$ws.Range("C10").Value = "Itga3"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 10.80562866666667
$ws.Range("H10").Value = 32.416886
$ws.Range("I10").Value = 0.2139607882293788
$ws.Range("J10").Value = 0.2139607882293788
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 5.608704333333333
$ws.Range("N10").Value = 16.826113
$ws.Range("O10").Value = 0.473454027939509
$ws.Range("P10").Value = 0.4734540279395091
$ws.Range("Q10").Value = 60.60557632712421
$ws.Range("R10").Value = 545.450186944118
$ws.Range("S10").Value = 0.1013005970083117
$ws.Range("T10").Value = 0.1013005970083117

Write-Host "Updated rows 2-10 with new LR-pair values"
